$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.128.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.561.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.001"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "289.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3797"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3290"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.44"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -9.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.141"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07376"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.42%  "
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("E13").Value = "  -3.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.828"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.899"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.565.15"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.31%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001095"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06658"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "85.81"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.466"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9997"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.16"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.72"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.142.53"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.270"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.98%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.557"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "151.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.868"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.739.45"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "121.41"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.123"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.71%  "
$ws.Range("E33").Value = "  -1.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.857"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.381"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.80%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08178"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.295"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.28%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06227"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.67%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02303"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2144"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.28%  "
$ws.Range("E41").Value = "  -4.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.08"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.000"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5997"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.79"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.65%  "
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5799"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.978"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.98%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "120.86"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.172"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06983"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.38%  "
